$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.952.92'
$ws.Range('E2').Value = '  -0.36%  '

$ws.Range('D3').Value = '2.513.40'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '''530.46'
$ws.Range('E5').Value = '  -1.84%  '

$ws.Range('D6').Value = '''138.83'
$ws.Range('E6').Value = '  -3.42%  '

$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.15%  '

$ws.Range('D8').Value = '''0.562'
$ws.Range('E8').Value = '  -1.77%  '

$ws.Range('D9').Value = '2.515.68'
$ws.Range('E9').Value = '  -0.44%  '

$ws.Range('D10').Value = '''0.0999'

$ws.Range('E11').Value = '  +0.59%  '

$ws.Range('D12').Value = '''5.45'
$ws.Range('E12').Value = '  -2.62%  '

$ws.Range('E13').Value = '  +0.18%  '

$ws.Range('D14').Value = '2.956.39'
$ws.Range('E14').Value = '  +0.38%  '

$ws.Range('D15').Value = '''22.98'
$ws.Range('E15').Value = '  -2.33%  '

$ws.Range('D16').Value = '58.895.06'
$ws.Range('E16').Value = '  -0.18%  '

$ws.Range('D18').Value = '2.509.36'
$ws.Range('E18').Value = '  -0.44%  '

$ws.Range('E19').Value = '  -1.98%  '

$ws.Range('E20').Value = '  -0.97%  '

$ws.Range('D21').Value = '''321.50'
$ws.Range('E21').Value = '  -1.08%  '

$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('E23').Value = '  +0.43%  '

$ws.Range('D24').Value = '''62.25'
$ws.Range('E24').Value = '  +0.10%  '

$ws.Range('E25').Value = '  -4.28%  '

$ws.Range('E26').Value = '  +1.90%  '

$ws.Range('E27').Value = '  +0.19%  '

$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('D29').Value = '''6.74'
$ws.Range('E29').Value = '  +0.98%  '

$ws.Range('D30').Value = '0.0₃0769'
$ws.Range('E30').Value = '  -0.84%  '

$ws.Range('E31').Value = '  -1.95%  '

$ws.Range('D32').Value = '''161.76'
$ws.Range('E32').Value = '  +3.39%  '

$ws.Range('E33').Value = '  +0.24%  '

$ws.Range('E34').Value = '  -6.93%  '

$ws.Range('E35').Value = '  -0.93%  '

$ws.Range('D36').Value = '''18.41'
$ws.Range('E36').Value = '  -1.37%  '

$ws.Range('E37').Value = '  -3.29%  '

$ws.Range('D38').Value = '''1.56'
$ws.Range('E38').Value = '  -1.93%  '

$ws.Range('D39').Value = '''36.95'
$ws.Range('E39').Value = '  +0.10%  '

$ws.Range('E40').Value = '  -1.85%  '

$ws.Range('E41').Value = '  -2.31%  '

$ws.Range('E42').Value = '  -9.32%  '

$ws.Range('D43').Value = '''278.66'
$ws.Range('E43').Value = '  -5.83%  '

$ws.Range('E44').Value = '  +0.20%  '

$ws.Range('E45').Value = '  +0.89%  '

$ws.Range('E46').Value = '  -0.94%  '

$ws.Range('D47').Value = '''0.0929'
$ws.Range('E47').Value = '  -0.06%  '

$ws.Range('D48').Value = '''121.42'
$ws.Range('E48').Value = '  -0.97%  '

$ws.Range('D49').Value = '''18.30'
$ws.Range('E49').Value = '  -1.48%  '

$ws.Range('D50').Value = '''0.0508'
$ws.Range('E50').Value = '  -1.04%  '

$ws.Range('E51').Value = '  -2.75%  '
